$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("F2").Value = 2.52
$ws.Range("G2").Value = 2.56
$ws.Range("H2").Value = 3.35
$ws.Range("J2").Value = 3.15
$ws.Range("K2").Value = 3.2
$ws.Range("T2").Value = 2
$ws.Range("AB2").Value = 8.800000000000001
$ws.Range("AD2").Value = 15

# Row 3 updates
$ws.Range("F3").Value = 1.45
$ws.Range("G3").Value = 1.46
$ws.Range("H3").Value = 9.800000000000001
$ws.Range("T3").Value = 2.44
$ws.Range("U3").Value = 1.67
$ws.Range("AA3").Value = 550
$ws.Range("AI3").Value = 220
$ws.Range("AM3").Value = 310

# Row 4 updates
$ws.Range("N4").Value = 4.4
$ws.Range("P4").Value = 2.16
$ws.Range("Q4").Value = 1.83
$ws.Range("R4").Value = 1.45
$ws.Range("T4").Value = 1.86
$ws.Range("U4").Value = 2.08
$ws.Range("X4").Value = 18
$ws.Range("AG4").Value = 9.800000000000001
$ws.Range("AH4").Value = 21
$ws.Range("AI4").Value = 85
$ws.Range("AO4").Value = 85
